$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the weekly cryptos data refresh (GitHub Actions scheduled update).
# Most cells are plain text updates; a handful of price cells in column D
# are plain decimal numbers that Excel would otherwise auto-convert to the
# Number type, so we briefly force a Text number format for those specific
# cells (matching the original inline-string/text cell type), then restore
# the default "General" number format.

$ws.Range("D2").Value = "63.677.51"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.320.64"
$ws.Range("E3").Value = "  +5.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.70"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.32"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.320.38"
$ws.Range("E8").Value = "  +5.32%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.72"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "3.868.23"
$ws.Range("E15").Value = "  +5.33%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "3.319.09"
$ws.Range("E17").Value = "  +5.35%  "
$ws.Range("D18").Value = "63.783.58"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.44"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.17"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  +4.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.16"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +5.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.72"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.96"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.17"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.08"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +7.57%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").Value = "  +3.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.70"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  +6.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0401"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "432.96"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").Value = "3.089.79"
$ws.Range("E41").Value = "  +5.39%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.35"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  +4.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +3.47%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.95"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +14.23%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.39"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.114"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.43%  "

Write-Host "Applied cryptos list update"
